$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the formatting of the last existing data row (385) down through
# the newly appended rows (386-464) so the new rows inherit the same
# cell styles (date format/border in column A, default elsewhere).
$srcRange = $ws.Range("A385:D385")
$dstRange = $ws.Range("A386:D464")
$srcRange.Copy()
$dstRange.PasteSpecial(-4122)

# New daily data: date (serial), nuovi pos., somma mobile 7gg.,
# somma mobile 7gg. per 100mila abitanti
$data = @(
    @(44460,1,5,156.2988433885589),
    @(44461,0,5,156.2988433885589),
    @(44462,1,5,156.2988433885589),
    @(44463,0,4,125.0390747108471),
    @(44464,0,4,125.0390747108471),
    @(44465,0,2,62.51953735542357),
    @(44466,0,2,62.51953735542357),
    @(44467,0,1,31.25976867771178),
    @(44468,0,1,31.25976867771178),
    @(44469,0,0,0),
    @(44470,0,0,0),
    @(44471,0,0,0),
    @(44472,1,1,31.25976867771178),
    @(44473,0,1,31.25976867771178),
    @(44474,1,2,62.51953735542357),
    @(44475,0,2,62.51953735542357),
    @(44476,1,3,93.77930603313536),
    @(44477,0,3,93.77930603313536),
    @(44478,0,3,93.77930603313536),
    @(44479,0,2,62.51953735542357),
    @(44480,0,2,62.51953735542357),
    @(44481,0,1,31.25976867771178),
    @(44482,0,1,31.25976867771178),
    @(44483,0,0,0),
    @(44484,0,0,0),
    @(44485,0,0,0),
    @(44486,0,0,0),
    @(44487,0,0,0),
    @(44488,0,0,0),
    @(44489,0,0,0),
    @(44490,0,0,0),
    @(44491,0,0,0),
    @(44492,0,0,0),
    @(44493,0,0,0),
    @(44494,0,0,0),
    @(44495,0,0,0),
    @(44496,0,0,0),
    @(44497,0,0,0),
    @(44498,0,0,0),
    @(44499,0,0,0),
    @(44500,0,0,0),
    @(44501,0,0,0),
    @(44502,0,0,0),
    @(44503,0,0,0),
    @(44504,0,0,0),
    @(44505,0,0,0),
    @(44506,0,0,0),
    @(44507,0,0,0),
    @(44508,0,0,0),
    @(44509,0,0,0),
    @(44510,0,0,0),
    @(44511,0,0,0),
    @(44512,0,0,0),
    @(44513,1,1,31.25976867771178),
    @(44514,0,1,31.25976867771178),
    @(44515,0,1,31.25976867771178),
    @(44516,2,3,93.77930603313536),
    @(44517,0,3,93.77930603313536),
    @(44518,0,3,93.77930603313536),
    @(44519,2,5,156.2988433885589),
    @(44520,0,4,125.0390747108471),
    @(44521,0,4,125.0390747108471),
    @(44522,2,6,187.5586120662707),
    @(44523,3,7,218.8183807439825),
    @(44524,2,9,281.3379180994061),
    @(44525,0,9,281.3379180994061),
    @(44526,1,8,250.0781494216943),
    @(44527,0,8,250.0781494216943),
    @(44528,0,8,250.0781494216943),
    @(44529,0,6,187.5586120662707),
    @(44530,0,3,93.77930603313536),
    @(44531,1,2,62.51953735542357),
    @(44532,0,2,62.51953735542357),
    @(44533,3,4,125.0390747108471),
    @(44534,0,4,125.0390747108471),
    @(44535,8,12,375.1172241325414),
    @(44536,0,12,375.1172241325414),
    @(44537,2,14,437.636761487965),
    @(44538,1,14,437.636761487965)
)

$startRow = 386
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
